$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 34: milestone changed from II to I, marked complete
$ws.Range("E34").Value = "I"
$ws.Range("F34").Value = "X"

# Row 35: milestone changed from II to I, marked complete
$ws.Range("E35").Value = "I"
$ws.Range("F35").Value = "X"

# Row 36: milestone changed from II to I, marked complete
$ws.Range("E36").Value = "I"
$ws.Range("F36").Value = "X"

# Update the sheet view to reflect scroll position / selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F34").Select()
